$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5 ("Salary") so it becomes row 6,
# and fill the freed row 5 with the new "income1" entry.
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "income1"
$ws.Range("B5").Value = 5000
$ws.Range("C5").Value = 45902.229537037034

# Copy the date formatting from the (now shifted) Salary row onto the two
# brand-new rows before writing their values, so C7/C8 pick up the same
# date style (s="1") as the rest of column C instead of "General".
$ws.Range("C6").Copy($ws.Range("C7:C8"))

# Append two more income rows at the end of the table.
$ws.Range("A7").Value = "income2"
$ws.Range("B7").Value = 3000
$ws.Range("C7").Value = 45889.229537037034

$ws.Range("A8").Value = "income3"
$ws.Range("B8").Value = 3500
$ws.Range("C8").Value = 45870.229537037034
